$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 17.10.2024"

# Row 6
$ws.Range("B6").Value = "18.10."
$ws.Range("C6").Value = "19.10."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 62327998"
$ws.Range("E6").Value = "41,97-"

# Row 7
$ws.Range("B7").Value = "19.10."
$ws.Range("C7").Value = "20.10."
$ws.Range("D7").Value = "PAYPAL TSXSPI"
$ws.Range("E7").Value = "45,36-"

# Row 8
$ws.Range("B8").Value = "22.10."
$ws.Range("C8").Value = "23.10."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,14-"

# Row 9
$ws.Range("B9").Value = "25.10."
$ws.Range("C9").Value = "26.10."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 59560162"
$ws.Range("E9").Value = "84,34-"

# Row 10 - cleared out (now an empty trailing row like row 11)
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
# E10's style index changes from the "amount" style (17) to the blank-row
# style (12) used by row 11's trailing E cell - copy formats across to match.
$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 27.10.2024"
$ws.Range("E12").Value = "196,81-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 01.11.2024"
